$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 34 with the new "Avances Etiquetado Roboflow" data (6/11/2025)
$ws.Range("D34").Value = (Get-Date -Year 2025 -Month 11 -Day 6).Date
$ws.Range("E34").Value = 75
$ws.Range("F34").Value = 443
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 1012
$ws.Range("J34").Value = "N/A"

# Move the active selection to F36, matching the state after entering the row
$ws.Range("F36").Select()
